# Trade #36 closed at 2026-02-18 00:14:43 - unknown UNKNOWN +0.000%
#
# This script applies a batch of live-trading-results updates:
#   1) Closes the open HighProbConvergence trade (#64) that was sitting in
#      "All Trades" row 65 / "HighProbConvergence" row 4.
#   2) Refreshes the rolled-up Summary + Strategy Status numbers to reflect
#      the newly closed trade.
#   3) Appends four freshly-opened trades (book-keeping entries #93-#96) to
#      "All Trades" and to each trade's own per-strategy sheet.
#
# NOTE: a handful of string values in this sheet happen to look like
# ISO dates ("2026-02-18"). Assigning such a literal to Range.Value makes
# Excel auto-convert it into a date serial (exactly like typing it into a
# live grid would). We force-preserve the literal text the same way a user
# would in the UI - a leading apostrophe - so the cell keeps its original
# plain-text/General-format shape instead of turning into a date.

$wb = $excel.ActiveWorkbook

function Set-Text {
    param($range, [string]$text)
    # Leading apostrophe = "treat as text" (mirrors typing '2026-02-18 into
    # a cell in the Excel UI); keeps date-shaped strings as literal text.
    $range.Value = "'" + $text
}

# ---------------------------------------------------------------------
# 1) Summary sheet - roll up totals after the close
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1499.81   # Current Capital
$wsSummary.Range("B4").Value = 0.92      # Total P&L $
$wsSummary.Range("B5").Value = 0.29      # Total P&L %
$wsSummary.Range("B6").Value = 64        # Total Trades
$wsSummary.Range("B8").Value = 26        # Losing Trades
$wsSummary.Range("B9").Value = 54.69     # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - HighProbConvergence row (row 3)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C3").Value = 100.06
$wsStatus.Range("D3").Value = 3
$wsStatus.Range("E3").Value = 0.06
$wsStatus.Range("F3").Value = 0.06
$wsStatus.Range("G3").Value = 66.67

# ---------------------------------------------------------------------
# 3) All Trades sheet - close trade #64 (row 65) + append rows 94-97
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# --- close out trade #64 (HighProbConvergence) ---
$wsAll.Cells.Item(65, 7).Value = 0.676563          # G65 Exit Price
$wsAll.Cells.Item(65, 8).Value = "CLOSED"          # H65 Status
$wsAll.Cells.Item(65, 9).Value = -23.1179          # I65 P&L %
$wsAll.Cells.Item(65, 10).Value = -0.2             # J65 P&L $
$wsAll.Cells.Item(65, 11).Value = 100.06            # K65 Capital After
$wsAll.Cells.Item(65, 12).Value = "early_exit"     # L65 Exit Reason
$wsAll.Cells.Item(65, 13).Value = 0.25             # M65 Duration (min)

# --- append newly opened trades ---
# Row 94: trade #93 (momentum)
Set-Text $wsAll.Cells.Item(94, 2) "2026-02-18"
$wsAll.Cells.Item(94, 1).Value = 93
$wsAll.Cells.Item(94, 3).Value = "00:14:35"
$wsAll.Cells.Item(94, 4).Value = "momentum"
$wsAll.Cells.Item(94, 5).Value = "DOWN"
$wsAll.Cells.Item(94, 6).Value = 0.88
$wsAll.Cells.Item(94, 8).Value = "OPEN"
$wsAll.Cells.Item(94, 9).Value = 0
$wsAll.Cells.Item(94, 10).Value = 0
$wsAll.Cells.Item(94, 11).Value = 100.0512903225807
$wsAll.Cells.Item(94, 13).Value = 0
$wsAll.Cells.Item(94, 14).Value = 0
$wsAll.Cells.Item(94, 15).Value = 0
$wsAll.Cells.Item(94, 16).Value = 0.9
$wsAll.Cells.Item(94, 17).Value = "Downward momentum: -40.404% over 10 samples"

# Row 95: trade #94 (HighProbConvergence)
$wsAll.Cells.Item(95, 1).Value = 94
Set-Text $wsAll.Cells.Item(95, 2) "2026-02-18"
$wsAll.Cells.Item(95, 3).Value = "00:14:36"
$wsAll.Cells.Item(95, 4).Value = "HighProbConvergence"
$wsAll.Cells.Item(95, 5).Value = "UP"
$wsAll.Cells.Item(95, 6).Value = 0.13
$wsAll.Cells.Item(95, 8).Value = "OPEN"
$wsAll.Cells.Item(95, 9).Value = 0
$wsAll.Cells.Item(95, 10).Value = 0
$wsAll.Cells.Item(95, 11).Value = 100.26
$wsAll.Cells.Item(95, 13).Value = 0
$wsAll.Cells.Item(95, 14).Value = 0
$wsAll.Cells.Item(95, 15).Value = 0
$wsAll.Cells.Item(95, 16).Value = 0.95
$wsAll.Cells.Item(95, 17).Value = "Mean reversion UP: price 37.89% below mean (z=-3.00)"

# Row 96: trade #95 (MarketMaking)
$wsAll.Cells.Item(96, 1).Value = 95
Set-Text $wsAll.Cells.Item(96, 2) "2026-02-18"
$wsAll.Cells.Item(96, 3).Value = "00:14:37"
$wsAll.Cells.Item(96, 4).Value = "MarketMaking"
$wsAll.Cells.Item(96, 5).Value = "DOWN"
$wsAll.Cells.Item(96, 6).Value = 0.86
$wsAll.Cells.Item(96, 8).Value = "OPEN"
$wsAll.Cells.Item(96, 9).Value = 0
$wsAll.Cells.Item(96, 10).Value = 0
$wsAll.Cells.Item(96, 11).Value = 99.58025471513901
$wsAll.Cells.Item(96, 13).Value = 0
$wsAll.Cells.Item(96, 14).Value = 0
$wsAll.Cells.Item(96, 15).Value = 0
$wsAll.Cells.Item(96, 16).Value = 0.65
$wsAll.Cells.Item(96, 17).Value = "Wide spread capture: 339 bps vs avg 210 bps"

# Row 97: trade #96 (EMAArbitrage)
$wsAll.Cells.Item(97, 1).Value = 96
Set-Text $wsAll.Cells.Item(97, 2) "2026-02-18"
$wsAll.Cells.Item(97, 3).Value = "00:14:37"
$wsAll.Cells.Item(97, 4).Value = "EMAArbitrage"
$wsAll.Cells.Item(97, 5).Value = "DOWN"
$wsAll.Cells.Item(97, 6).Value = 0.83
$wsAll.Cells.Item(97, 8).Value = "OPEN"
$wsAll.Cells.Item(97, 9).Value = 0
$wsAll.Cells.Item(97, 10).Value = 0
$wsAll.Cells.Item(97, 11).Value = 100.1258137286497
$wsAll.Cells.Item(97, 13).Value = 0
$wsAll.Cells.Item(97, 14).Value = 0
$wsAll.Cells.Item(97, 15).Value = 0
$wsAll.Cells.Item(97, 16).Value = 0.9
$wsAll.Cells.Item(97, 17).Value = "EMA:down, RSI:0.0, ROC:-40.40% | 2/3 DOWN"

# ---------------------------------------------------------------------
# 4) momentum sheet - append row 21 (trade #93)
# ---------------------------------------------------------------------
$wsMomentum = $wb.Worksheets.Item("momentum")
$wsMomentum.Cells.Item(21, 1).Value = 93
Set-Text $wsMomentum.Cells.Item(21, 2) "2026-02-18"
$wsMomentum.Cells.Item(21, 3).Value = "00:14:35"
$wsMomentum.Cells.Item(21, 4).Value = "momentum"
$wsMomentum.Cells.Item(21, 5).Value = "DOWN"
$wsMomentum.Cells.Item(21, 6).Value = 0.88
$wsMomentum.Cells.Item(21, 8).Value = "OPEN"
$wsMomentum.Cells.Item(21, 9).Value = 0
$wsMomentum.Cells.Item(21, 10).Value = 0
$wsMomentum.Cells.Item(21, 11).Value = 100.0512903225807
$wsMomentum.Cells.Item(21, 12).Value = 0
$wsMomentum.Cells.Item(21, 13).Value = 0
$wsMomentum.Cells.Item(21, 14).Value = 0.9
$wsMomentum.Cells.Item(21, 15).Value = "Downward momentum: -40.404% over 10 samples"
$wsMomentum.Cells.Item(21, 17).Value = 0

# ---------------------------------------------------------------------
# 5) HighProbConvergence sheet - close trade #64 (row 4) + append row 10
# ---------------------------------------------------------------------
$wsHPC = $wb.Worksheets.Item("HighProbConvergence")

# --- close out trade #64 ---
$wsHPC.Cells.Item(4, 7).Value = 0.676563
$wsHPC.Cells.Item(4, 8).Value = "CLOSED"
$wsHPC.Cells.Item(4, 9).Value = -23.1179
$wsHPC.Cells.Item(4, 10).Value = -0.2
$wsHPC.Cells.Item(4, 11).Value = 100.06
$wsHPC.Cells.Item(4, 16).Value = "early_exit"
$wsHPC.Cells.Item(4, 17).Value = 0.25

# --- append trade #94 ---
$wsHPC.Cells.Item(10, 1).Value = 94
Set-Text $wsHPC.Cells.Item(10, 2) "2026-02-18"
$wsHPC.Cells.Item(10, 3).Value = "00:14:36"
$wsHPC.Cells.Item(10, 4).Value = "HighProbConvergence"
$wsHPC.Cells.Item(10, 5).Value = "UP"
$wsHPC.Cells.Item(10, 6).Value = 0.13
$wsHPC.Cells.Item(10, 8).Value = "OPEN"
$wsHPC.Cells.Item(10, 9).Value = 0
$wsHPC.Cells.Item(10, 10).Value = 0
$wsHPC.Cells.Item(10, 11).Value = 100.26
$wsHPC.Cells.Item(10, 12).Value = 0
$wsHPC.Cells.Item(10, 13).Value = 0
$wsHPC.Cells.Item(10, 14).Value = 0.95
$wsHPC.Cells.Item(10, 15).Value = "Mean reversion UP: price 37.89% below mean (z=-3.00)"
$wsHPC.Cells.Item(10, 17).Value = 0

# ---------------------------------------------------------------------
# 6) MarketMaking sheet - append row 32 (trade #95)
# ---------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Cells.Item(32, 1).Value = 95
Set-Text $wsMM.Cells.Item(32, 2) "2026-02-18"
$wsMM.Cells.Item(32, 3).Value = "00:14:37"
$wsMM.Cells.Item(32, 4).Value = "MarketMaking"
$wsMM.Cells.Item(32, 5).Value = "DOWN"
$wsMM.Cells.Item(32, 6).Value = 0.86
$wsMM.Cells.Item(32, 8).Value = "OPEN"
$wsMM.Cells.Item(32, 9).Value = 0
$wsMM.Cells.Item(32, 10).Value = 0
$wsMM.Cells.Item(32, 11).Value = 99.58025471513901
$wsMM.Cells.Item(32, 12).Value = 0
$wsMM.Cells.Item(32, 13).Value = 0
$wsMM.Cells.Item(32, 14).Value = 0.65
$wsMM.Cells.Item(32, 15).Value = "Wide spread capture: 339 bps vs avg 210 bps"
$wsMM.Cells.Item(32, 17).Value = 0

# ---------------------------------------------------------------------
# 7) EMAArbitrage sheet - append row 8 (trade #96)
# ---------------------------------------------------------------------
$wsEMA = $wb.Worksheets.Item("EMAArbitrage")
$wsEMA.Cells.Item(8, 1).Value = 96
Set-Text $wsEMA.Cells.Item(8, 2) "2026-02-18"
$wsEMA.Cells.Item(8, 3).Value = "00:14:37"
$wsEMA.Cells.Item(8, 4).Value = "EMAArbitrage"
$wsEMA.Cells.Item(8, 5).Value = "DOWN"
$wsEMA.Cells.Item(8, 6).Value = 0.83
$wsEMA.Cells.Item(8, 8).Value = "OPEN"
$wsEMA.Cells.Item(8, 9).Value = 0
$wsEMA.Cells.Item(8, 10).Value = 0
$wsEMA.Cells.Item(8, 11).Value = 100.1258137286497
$wsEMA.Cells.Item(8, 12).Value = 0
$wsEMA.Cells.Item(8, 13).Value = 0
$wsEMA.Cells.Item(8, 14).Value = 0.9
$wsEMA.Cells.Item(8, 15).Value = "EMA:down, RSI:0.0, ROC:-40.40% | 2/3 DOWN"
$wsEMA.Cells.Item(8, 17).Value = 0
